$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A12").Value = "27.3.2019"
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = "profile-komponentti, refaktorointia"

$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").VerticalAlignment = -4160
$ws.Range("B12").VerticalAlignment = -4160
$ws.Range("C12").VerticalAlignment = -4160
$ws.Range("C12").WrapText = $true
